$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh (prices + 1h volume %) pulled on Mon Jul 15 23:10:19 UTC 2024.
# Rows 31-33 and 50-51 also change rank order (coin name/link swap).

# Row 2
$ws.Range("D2").Value = "64.345.89"
$ws.Range("E2").Value = "  +5.64%  "

# Row 3
$ws.Range("D3").Value = "3.474.61"
$ws.Range("E3").Value = "  +6.94%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").Value = "`'582.43"
$ws.Range("E5").Value = "  +6.84%  "

# Row 6
$ws.Range("E6").Value = "  +7.49%  "

# Row 7
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("D8").Value = "3.481.39"
$ws.Range("E8").Value = "  +6.96%  "

# Row 9
$ws.Range("D9").Value = "`'0.534"
$ws.Range("E9").Value = "  +1.24%  "

# Row 11
$ws.Range("E11").Value = "  +7.28%  "

# Row 12
$ws.Range("D12").Value = "`'0.441"
$ws.Range("E12").Value = "  +1.79%  "

# Row 13
$ws.Range("D13").Value = "4.064.78"
$ws.Range("E13").Value = "  +6.74%  "

# Row 14
$ws.Range("E14").Value = "  -0.50%  "

# Row 15
$ws.Range("E15").Value = "  +7.75%  "

# Row 16
$ws.Range("D16").Value = "`'27.65"
$ws.Range("E16").Value = "  +4.74%  "

# Row 17
$ws.Range("D17").Value = "64.394.14"
$ws.Range("E17").Value = "  +5.76%  "

# Row 18
$ws.Range("D18").Value = "3.465.21"
$ws.Range("E18").Value = "  +6.61%  "

# Row 19
$ws.Range("D19").Value = "`'6.44"
$ws.Range("E19").Value = "  +2.09%  "

# Row 20
$ws.Range("D20").Value = "`'14.36"
$ws.Range("E20").Value = "  +6.76%  "

# Row 21
$ws.Range("D21").Value = "`'397.23"
$ws.Range("E21").Value = "  +4.80%  "

# Row 22
$ws.Range("D22").Value = "`'8.52"
$ws.Range("E22").Value = "  +0.85%  "

# Row 23
$ws.Range("D23").Value = "`'0.545"
$ws.Range("E23").Value = "  +2.24%  "

# Row 24
$ws.Range("D24").Value = "`'0.996"
$ws.Range("E24").Value = "  -0.46%  "

# Row 25
$ws.Range("D25").Value = "`'72.04"
$ws.Range("E25").Value = "  +2.83%  "

# Row 26
$ws.Range("D26").Value = "`'0.0000110"
$ws.Range("E26").Value = "  +20.37%  "

# Row 27
$ws.Range("D27").Value = "`'9.48"
$ws.Range("E27").Value = "  +9.80%  "

# Row 28
$ws.Range("E28").Value = "  +6.56%  "

# Row 29
$ws.Range("E29").Value = "  -0.13%  "

# Row 30
$ws.Range("E30").Value = "  +13.58%  "

# Row 31
$ws.Range("B31").Value = "RenderToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D31").Value = "`'6.73"
$ws.Range("E31").Value = "  +8.42%  "

# Row 32
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "`'2.06"
$ws.Range("E32").Value = "  +6.11%  "

# Row 33
$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").Value = "`'5.85"
$ws.Range("E33").Value = "  +7.93%  "

# Row 34
$ws.Range("D34").Value = "`'23.86"
$ws.Range("E34").Value = "  +5.43%  "

# Row 35
$ws.Range("E35").Value = "  -0.11%  "

# Row 36
$ws.Range("D36").Value = "`'6.93"
$ws.Range("E36").Value = "  +4.33%  "

# Row 37
$ws.Range("E37").Value = "  +5.40%  "

# Row 38
$ws.Range("D38").Value = "`'158.98"
$ws.Range("E38").Value = "  -0.12%  "

# Row 39
$ws.Range("D39").Value = "`'28.39"
$ws.Range("E39").Value = "  +7.64%  "

# Row 40
$ws.Range("D40").Value = "`'0.0787"
$ws.Range("E40").Value = "  +8.91%  "

# Row 41
$ws.Range("D41").Value = "`'1.88"
$ws.Range("E41").Value = "  +9.37%  "

# Row 42
$ws.Range("D42").Value = "2.876.76"
$ws.Range("E42").Value = "  +2.67%  "

# Row 43
$ws.Range("E43").Value = "  +3.29%  "

# Row 44
$ws.Range("D44").Value = "`'0.785"
$ws.Range("E44").Value = "  +7.24%  "

# Row 45
$ws.Range("D45").Value = "`'4.43"
$ws.Range("E45").Value = "  +3.37%  "

# Row 46
$ws.Range("D46").Value = "`'42.13"
$ws.Range("E46").Value = "  +5.09%  "

# Row 47
$ws.Range("E47").Value = "  +10.41%  "

# Row 48
$ws.Range("D48").Value = "3.515.22"
$ws.Range("E48").Value = "  +6.82%  "

# Row 49
$ws.Range("D49").Value = "`'22.83"
$ws.Range("E49").Value = "  +5.82%  "

# Row 50
$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").Value = "`'2.12"
$ws.Range("E50").Value = "  +22.31%  "

# Row 51
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").Value = "`'6.42"
$ws.Range("E51").Value = "  +2.76%  "
